$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "H2" = 84
    "I2" = 229
    "J2" = 851
    "K2" = 2
    "L2" = 214
    "M2" = 7
    "N2" = 140
    "O2" = 0
    "P2" = 3
    "Q2" = 3
    "R2" = 3
    "S2" = 85
    "T2" = 167
    "U2" = 10
    "V2" = 1286
    "W2" = 1
    "X2" = 1330
    "Y2" = 2
    "Z2" = 20
    "AA2" = 15
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
